$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price (D) column updates: force text entry (values are plain text in the
# source data, e.g. thousands-dot formatted numbers) so Excel does not
# reinterpret/reformat them as numeric values.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.865.52"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  -2.24%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.778.44"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -2.85%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.008"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  +0.59%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "307.31"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -1.70%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4225"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -1.62%  "

$ws.Range("E8").Value = "  -1.06%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07169"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -1.66%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8378"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -3.25%  "

$ws.Range("E11").Value = "  -1.94%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.753.49"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -5.93%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.250"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -2.82%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.341"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -2.75%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.06812"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -1.67%  "

$ws.Range("E16").Value = "  +0.73%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "79.11"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -1.66%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008678"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -2.73%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.007"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +0.54%  "

$ws.Range("E20").Value = "  -2.98%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "26.725.22"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -3.07%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.011"
$ws.Range("D22").ClearFormats()

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "11.07"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +2.23%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.981.14"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -3.56%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.922"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -3.02%  "

$ws.Range("E26").Value = "  -0.66%  "

$ws.Range("E27").Value = "  -4.38%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "5.043"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -1.73%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "114.35"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -0.03%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.625"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -11.63%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08949"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +0.89%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.7197"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -4.64%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.836"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -4.82%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.323"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -4.79%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.090"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -4.06%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.007"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +0.58%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.081"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -0.78%  "

$ws.Range("E38").Value = "  -2.41%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05082"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -4.60%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.4916"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -3.37%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.519"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -10.05%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.080"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -7.52%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "7.920"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -5.01%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.008"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +0.69%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "104.44"
$ws.Range("D46").ClearFormats()

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.08"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -3.68%  "

$ws.Range("E48").Value = "  -4.20%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.4477"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -4.42%  "

$ws.Range("E50").Value = "  -2.84%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.718"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -1.06%  "
